$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "leonardo "
$ws.Range("A4").Value = "sara"

$ws.Range("A4").Select()
